$d = $word.ActiveDocument

# The run "Cen, Sur\" / "endra" / " " was a stray mid-word line break
# (with spell-check markers around "endra"); merge it back into a single
# run reading "Cen, Surendra ".
$d.Content.Find.Execute("Sur\endra ", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Surendra ", 2) | Out-Null

# Find the "Members: ... Michelle ..." paragraph (the one team roster that
# is still missing the separator rule the other "Members:" paragraphs use)
# and give it the same single-line bottom border.
$target = $null
foreach ($para in $d.Paragraphs) {
    if ($para.Range.Text -like "Members: Michelle*") {
        $target = $para
        break
    }
}

if ($target -ne $null) {
    $target.Format.Borders.Item(-3).LineStyle = 1
    $target.Format.Borders.Item(-3).LineWidth = 2
    $target.Format.Borders.DistanceFromBottom = 1
    $target.Format.Borders.Item(-3).ColorIndex = 0
}
